$d = $word.ActiveDocument

# Locate the existing bullet "Explore various interpretability methods "
$rng = $d.Content
$found = $rng.Find.Execute("Explore various interpretability methods ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans the found text; split a new paragraph off right after
    # it, inheriting the same list/paragraph/run formatting.
    $rng.InsertParagraphAfter()

    # The newly created (empty) paragraph is the one right after the
    # paragraph containing the found text.
    $newPara = $rng.Paragraphs(1).Next()
    $newRange = $newPara.Range
    $newRange.Text = "What is the difference between interpretability and explainability "
}
